# Fruta / hortaliza, semanal
# Rotate the weekly price-report values (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) across the
# 5 data rows (2-6): each row takes on the values that used to sit two rows
# below it, wrapping around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..6 (columns D, M, N, O, P, R, S), per the target data.
$rowData = @{
    2 = @{ D = 45106; M = 50;  N = 10000; O = 10000; P = 10000; R = "Región de O'Higgins"; S = 556 }
    3 = @{ D = 44719; M = 50;  N = 20000; O = 21000; P = 20400; R = "Provincia de Limarí";  S = 1133 }
    4 = @{ D = 44362; M = 100; N = 19000; O = 20000; P = 19500; R = "Provincia de Curicó";  S = 1083 }
    5 = @{ D = 45084; M = 100; N = 17000; O = 18000; P = 17500; R = "Región de O'Higgins"; S = 972 }
    6 = @{ D = 44320; M = 50;  N = 18000; O = 20000; P = 18800; R = "Provincia de Limarí";  S = 1044 }
}

foreach ($r in 2..6) {
    $vals = $rowData[$r]

    $ws.Cells.Item($r, 4).Value2 = $vals.D    # Fecha
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # Volumen
    $ws.Cells.Item($r, 14).Value2 = $vals.N   # Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $vals.O   # Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value2 = $vals.R   # Origen
    $ws.Cells.Item($r, 19).Value2 = $vals.S   # Precio $/Kg
}
